# Fund unit setting workbook: add a "Carry% (>4x)" column (G) with values,
# give it its own named cell style ("Normal 2" / Arial 11), resize the
# header row and adjust row heights / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New named cell style used for column G (header + data) ---------------
$style = $wb.Styles.Add("Normal 2")
$style.Font.Name = "Arial"
$style.Font.Size = 11

# --- Header row ------------------------------------------------------------
$ws.Range("G1").Value = "Carry% (>4x)"
$ws.Range("G1:G3").Style = "Normal 2"

# --- Data rows ---------------------------------------------------------------
$ws.Range("G2").Value = 30
$ws.Range("G3").Value = 30

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 14.25
$ws.Rows.Item(2).RowHeight = 13.5
$ws.Rows.Item(3).RowHeight = 13.5

# --- Selection ---------------------------------------------------------------
$ws.Range("G1:G3").Select()
